$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "20.9.2025"
$ws.Range("B11").Value = 0.5
$ws.Range("C11").Value = 0.625
$ws.Range("D11").Value = 0.83333333333333337
$ws.Range("E11").Value = 0.875

$ws.Range("B11:E11").NumberFormat = "h:mm AM/PM"

$ws.Range("F13").Select()
